$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 55: No=54, Name=Юрій, Action=Ввійшов, Room=0, Time=..., Status=in
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "Юрій"
$ws.Cells.Item(55, 3).Value = "Ввійшов"
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 5).Value = 43953.378890231485
$ws.Range("E54").Copy()
$ws.Cells.Item(55, 5).PasteSpecial(-4122)
$ws.Cells.Item(55, 6).Value = "in"

# Row 56: No=55, Name=Юрій, Action=Вийшов, Room=0, Time=..., Status=out
$ws.Cells.Item(56, 1).Value = 55
$ws.Cells.Item(56, 2).Value = "Юрій"
$ws.Cells.Item(56, 3).Value = "Вийшов"
$ws.Cells.Item(56, 4).Value = 0
$ws.Cells.Item(56, 5).Value = 43953.37916302083
$ws.Range("E54").Copy()
$ws.Cells.Item(56, 5).PasteSpecial(-4122)
$ws.Cells.Item(56, 6).Value = "out"
